$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27 (A27 = "cholesterol HDL", a duplicate of A15) is removed; subsequent
# rows shift up by one.
$ws.Rows.Item(27).Delete()

# Update the view to match: scrolled so row 19 is the top row, with E49 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("E49").Select()
